# "Generate Report for Archive"
# Flip every "Ready for handoff" status to "In Translation" on all three
# sheets, then shrink the now-narrower Status columns to fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Status columns are E (zh-cn) and F (de-de), rows 2-4
$wsOverview.Range("E2:E4").Value = "In Translation"
$wsOverview.Range("F2:F4").Value = "In Translation"

# Per-locale sheets: Status is column C, rows 2-4
$wsZhCn.Range("C2:C4").Value = "In Translation"
$wsDeDe.Range("C2:C4").Value = "In Translation"

# Re-fit the Status columns now that the text is shorter
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
